$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-19 Wednesday" "2025-02-20 Thursday"

Replace-Text "79÷7=11, 2" "32÷4=8, 0"
Replace-Text "82÷4=20, 2" "95÷7=13, 4"
Replace-Text "56÷9=6, 2" "58÷9=6, 4"
Replace-Text "84÷7=12, 0" "95÷3=31, 2"
Replace-Text "94÷5=18, 4" "45÷6=7, 3"

Replace-Text "57÷7=8, 1" "78÷9=8, 6"
Replace-Text "55÷3=18, 1" "31÷9=3, 4"
Replace-Text "12÷7=1, 5" "72÷5=14, 2"
Replace-Text "52÷4=13, 0" "40÷3=13, 1"
Replace-Text "43÷4=10, 3" "34÷9=3, 7"

Replace-Text "24÷9=2, 6" "50÷5=10, 0"
Replace-Text "29÷8=3, 5" "21÷4=5, 1"
Replace-Text "20÷8=2, 4" "49÷4=12, 1"
Replace-Text "73÷7=10, 3" "15÷5=3, 0"
Replace-Text "57÷4=14, 1" "79÷8=9, 7"

Replace-Text "81÷5=16, 1" "60÷4=15, 0"
Replace-Text "96÷6=16, 0" "75÷6=12, 3"
Replace-Text "78÷2=39, 0" "36÷4=9, 0"
Replace-Text "89÷6=14, 5" "96÷2=48, 0"
Replace-Text "15÷4=3, 3" "85÷8=10, 5"

Replace-Text "19÷4=4, 3" "43÷2=21, 1"
Replace-Text "64÷6=10, 4" "80÷7=11, 3"
Replace-Text "88÷2=44, 0" "22÷8=2, 6"
Replace-Text "66÷4=16, 2" "56÷2=28, 0"
Replace-Text "26÷9=2, 8" "38÷6=6, 2"
